# Update the "login" sheet: columns G, H (username) and I (email) for rows 2-21
# contain literal strings like "EthanBaker24" / "EthanBaker24@gmail.com" that need
# their trailing "24" year-suffix changed to "25".
$wb = $excel.ActiveWorkbook
$loginWs = $wb.Worksheets.Item("login")

for ($r = 2; $r -le 21; $r++) {
    $gCell = $loginWs.Cells.Item($r, 7)   # column G
    $hCell = $loginWs.Cells.Item($r, 8)   # column H
    $iCell = $loginWs.Cells.Item($r, 9)   # column I

    $gVal = $gCell.Value2
    $hVal = $hCell.Value2
    $iVal = $iCell.Value2

    if ($gVal -ne $null -and $gVal -like "*24") {
        $gCell.Value = ($gVal -replace '24$', '25')
    }
    if ($hVal -ne $null -and $hVal -like "*24") {
        $hCell.Value = ($hVal -replace '24$', '25')
    }
    if ($iVal -ne $null -and $iVal -like "*24@gmail.com") {
        $iCell.Value = ($iVal -replace '24@gmail\.com$', '25@gmail.com')
    }
}

# Update the "Sheet1" sheet: cell I2 holds the numeric year value (24 -> 25).
# Columns A, B, C use CONCATENATE formulas referencing $I$2 and will recalc
# automatically once the value changes.
$dataWs = $wb.Worksheets.Item("Sheet1")
$dataWs.Range("I2").Value = 25
